# Update cryptos list with latest scraped values (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Price($row, $value) {
    # Force text so values like "379.98" / "1.00" / "7.90" keep their exact
    # literal form instead of being coerced to a number, while leaving the
    # cell style unchanged afterwards.
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

function Set-Vol($row, $value) {
    $ws.Range("E$row").Value = "  $value  "
}

# Row 2 - Bitcoin
Set-Price 2 "51.046.67"
Set-Vol 2 "+0.16%"

# Row 3 - Ethereum
Set-Price 3 "2.955.66"
Set-Vol 3 "+0.84%"

# Row 5 - BNB
Set-Price 5 "379.98"
Set-Vol 5 "+1.61%"

# Row 6 - Solana
Set-Price 6 "102.16"
Set-Vol 6 "+0.82%"

# Row 7 - XRP
Set-Price 7 "0.544"
Set-Vol 7 "+1.83%"

# Row 8 - USDC
Set-Vol 8 "-0.03%"

# Row 9 - Cardano
Set-Price 9 "0.587"
Set-Vol 9 "+1.31%"

# Row 10 - Avalanche
Set-Price 10 "36.61"
Set-Vol 10 "+0.85%"

# Row 11 - TRON
Set-Price 11 "0.138"
Set-Vol 11 "-1.03%"

# Row 12 - Dogecoin
Set-Price 12 "0.0852"
Set-Vol 12 "+2.21%"

# Row 13 - WrappedliquidstakedEther2.0
Set-Price 13 "3.418.62"
Set-Vol 13 "+0.80%"

# Row 14 - now Polkadot (was Chainlink)
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-Price 14 "7.75"
Set-Vol 14 "+6.02%"

# Row 15 - now Chainlink (was Polkadot)
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-Price 15 "18.28"
Set-Vol 15 "+2.15%"

# Row 16 - Uniswap
Set-Price 16 "11.93"
Set-Vol 16 "+67.32%"

# Row 17 - WrappedEther
Set-Price 17 "2.961.32"
Set-Vol 17 "+0.81%"

# Row 18 - Polygon
Set-Price 18 "0.999"
Set-Vol 18 "+2.39%"

# Row 19 - WrappedBTC
Set-Price 19 "51.102.98"
Set-Vol 19 "+0.31%"

# Row 20 - ImmutableX
Set-Price 20 "3.09"
Set-Vol 20 "-1.38%"

# Row 21 - InternetComputer(DFINITY)
Set-Price 21 "12.36"
Set-Vol 21 "-1.09%"

# Row 22 - ShibaInu
Set-Price 22 "0.0₃0961"
Set-Vol 22 "+0.84%"

# Row 23 - Litecoin
Set-Price 23 "69.95"
Set-Vol 23 "+2.55%"

# Row 24 - PancakeSwap
Set-Price 24 "3.28"
Set-Vol 24 "+14.14%"

# Row 25 - BitcoinCash
Set-Price 25 "267.24"
Set-Vol 25 "+1.18%"

# Row 26 - Filecoin
Set-Price 26 "7.90"
Set-Vol 26 "-1.18%"

# Row 27 - now Dai (was RenderToken)
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-Price 27 "1.00"
Set-Vol 27 "-0.03%"

# Row 28 - now RenderToken (was Dai)
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-Price 28 "7.16"
Set-Vol 28 "-6.98%"

# Row 29 - Kaspa
Set-Price 29 "0.166"
Set-Vol 29 "-0.79%"

# Row 30 - EthereumClassic
Set-Vol 30 "+1.15%"

# Row 31 - Hedera
Set-Vol 31 "-2.17%"

# Row 32 - Cosmos
Set-Price 32 "10.39"
Set-Vol 32 "+5.88%"

# Row 33 - InjectiveProtocol
Set-Price 33 "34.27"
Set-Vol 33 "+2.30%"

# Row 34 - OKB
Set-Price 34 "50.97"
Set-Vol 34 "+0.33%"

# Row 35 - Toncoin
Set-Price 35 "2.07"
Set-Vol 35 "+2.44%"

# Row 36 - VeChain
Set-Price 36 "0.0435"
Set-Vol 36 "-3.23%"

# Row 37 - FirstDigitalUSD
Set-Vol 37 "-0.01%"

# Row 38 - LidoDAOToken
Set-Price 38 "3.25"
Set-Vol 38 "+9.70%"

# Row 39 - Stellar
Set-Vol 39 "+2.04%"

# Row 40 - now Celestia (was ARBITRUM)
$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-Price 40 "16.56"
Set-Vol 40 "+1.46%"

# Row 41 - now ARBITRUM (was Celestia)
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-Price 41 "1.83"
Set-Vol 41 "+3.15%"

# Row 42 - Monero
Set-Price 42 "124.52"
Set-Vol 42 "+3.42%"

# Row 43 - Stacks
Set-Vol 43 "-1.42%"

# Row 44 - EnergySwap
Set-Price 44 "21.57"
Set-Vol 44 "+3.33%"

# Row 45 - NEARProtocol
Set-Price 45 "3.54"
Set-Vol 45 "+10.45%"

# Row 46 - ApeXProtocol
Set-Price 46 "2.40"
Set-Vol 46 "+3.94%"

# Row 47 - WEMIXToken
Set-Vol 47 "-1.05%"

# Row 48 - Maker
Set-Price 48 "2.049.44"
Set-Vol 48 "+4.19%"

# Row 49 - TheGraph
Set-Price 49 "0.267"
Set-Vol 49 "-4.13%"

# Row 50 - BEAM
Set-Price 50 "0.0319"
Set-Vol 50 "-6.52%"

# Row 51 - THORChain
Set-Price 51 "5.40"
Set-Vol 51 "+7.57%"
